# Update build timestamp in the "mines - January 30" release version strings.
$wb = $excel.ActiveWorkbook

$newStamp = "February 02 2026 12.49.33 EST"

$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("A2").Value = "Version: mines - January 30 (built on $newStamp)"
$aboutWs.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Xiegou Coal Mine, China, M3823, version 'mines - January 30 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 9; $r++) {
    $cell = $dataWs.Cells.Item($r, 19)  # Column S = 19
    $cell.Value = "mines - January 30 (built on $newStamp)"
}
